# Update p-values in the Fig3 table (commit: "no RECOVERY, Status instead")
$d = $word.ActiveDocument

# Each pair is (old value, new value). Values are unique across the
# document's text runs, so a simple whole-document Find/Replace (with
# MatchWholeWord) targets exactly the intended table cell each time.
# Order matters: the "0.98" -> "0.99" replacement must run before the
# "0.90" -> "0.98" replacement creates a fresh "0.98" elsewhere in the
# table, so we process them in this sequence.
$changes = @(
    @("0.82", "0.77"),
    @("0.69", "0.60"),
    @("0.36", "0.93"),
    @("0.63", "0.67"),
    @("0.98", "0.99"),
    @("0.94", "0.96"),
    @("0.11", "0.19"),
    @("0.90", "0.98"),
    @("0.02", "0.05")
)

foreach ($change in $changes) {
    $old = $change[0]
    $new = $change[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
